$wb = $excel.ActiveWorkbook

# ----- Sheet 1: "Data" -----
$ws1 = $wb.Worksheets.Item(1)

# LoginTest block: row 4's Runmode flips Y -> N
$ws1.Range("A4").Value = "N"

# TestB block is rebuilt into the new ProfileTest block (header + one real
# profile row instead of the old 4 dummy rows / 7 columns)
$ws1.Range("A7").Value = "ProfileTest"
$ws1.Range("D8").Value = "ProfileName"
$ws1.Range("E8:G8").ClearContents()

$ws1.Range("B9").Value = "vaibhavcool20@protonmail.com"
$ws1.Range("D9").Value = "Vaibhav Gupta"
$ws1.Range("E9:G9").ClearContents()
$ws1.Range("C9").Value = "xxxxxxxx"

# Password placeholder on the LoginTest row gains an extra "x" (edited last
# so it keeps its shared-string slot at the tail of the table)
$ws1.Range("C3").Value = "xxxxxxxx"

# The old TestB table had 3 extra data rows that no longer exist; removing
# them shifts the TestC block up directly under the new ProfileTest block.
$ws1.Range("A10:A12").EntireRow.Delete()

$ws1.Columns.Item(4).ColumnWidth = 13.14

# ----- Sheet 2: "Testcase" -----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A3").Value = "ProfileTest"
$ws2.Range("A3").Select()

$ws1.Range("C3").Select()
